# Auto-generated edit script: updates FFXIV leve-profit numeric columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# per refreshed Universalis market data, across 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 866964.75
$ws.Range("I11").Value = 866964.75
$ws.Range("K11").Value = 866964.75
$ws.Range("M11").Value = -866824.75
$ws.Range("H15").Value = 251.14
$ws.Range("I15").Value = 251.14
$ws.Range("K15").Value = 753.42
$ws.Range("M15").Value = -584.42
$ws.Range("H129").Value = 787.7692
$ws.Range("J129").Value = 910.8570999999999
$ws.Range("L129").Value = 2732.5713
$ws.Range("N129").Value = -12732.5713
$ws.Range("H132").Value = 2452204
$ws.Range("I132").Value = 2779037.8
$ws.Range("J132").Value = 951
$ws.Range("K132").Value = 8337113.399999999
$ws.Range("L132").Value = 2853
$ws.Range("M132").Value = -8334583.399999999
$ws.Range("N132").Value = -7913

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2266.5918
$ws.Range("I61").Value = 1285
$ws.Range("J61").Value = 3957.111
$ws.Range("K61").Value = 1285
$ws.Range("L61").Value = 3957.111
$ws.Range("M61").Value = -1073
$ws.Range("N61").Value = -4381.111
$ws.Range("H136").Value = 2266.5918
$ws.Range("I136").Value = 1285
$ws.Range("J136").Value = 3957.111
$ws.Range("K136").Value = 3855
$ws.Range("L136").Value = 11871.333
$ws.Range("M136").Value = -1305
$ws.Range("N136").Value = -16971.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1881.5186
$ws.Range("I20").Value = 1621.1666
$ws.Range("J20").Value = 2089.8
$ws.Range("K20").Value = 1621.1666
$ws.Range("L20").Value = 2089.8
$ws.Range("M20").Value = -1374.1666
$ws.Range("N20").Value = -2583.8
$ws.Range("H22").Value = 268.25
$ws.Range("I22").Value = 310.5
$ws.Range("J22").Value = 226
$ws.Range("K22").Value = 310.5
$ws.Range("L22").Value = 226
$ws.Range("M22").Value = -137.5
$ws.Range("N22").Value = -572
$ws.Range("H107").Value = 1333.3334
$ws.Range("I107").Value = 1333.3334
$ws.Range("K107").Value = 1333.3334
$ws.Range("M107").Value = 586.6666
$ws.Range("H134").Value = 1255.3611
$ws.Range("I134").Value = 1089.1666
$ws.Range("J134").Value = 2086.3333
$ws.Range("K134").Value = 3267.4998
$ws.Range("L134").Value = 6258.999899999999
$ws.Range("M134").Value = -732.4998000000001
$ws.Range("N134").Value = -11328.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3658.7778
$ws.Range("I31").Value = 4155.6
$ws.Range("J31").Value = 3410.3667
$ws.Range("K31").Value = 4155.6
$ws.Range("L31").Value = 3410.3667
$ws.Range("M31").Value = -3860.6
$ws.Range("N31").Value = -4000.3667
$ws.Range("H34").Value = 3658.7778
$ws.Range("I34").Value = 4155.6
$ws.Range("J34").Value = 3410.3667
$ws.Range("K34").Value = 4155.6
$ws.Range("L34").Value = 3410.3667
$ws.Range("M34").Value = -3953.6
$ws.Range("N34").Value = -3814.3667
$ws.Range("H58").Value = 1714.7587
$ws.Range("I58").Value = 1638.0741
$ws.Range("J58").Value = 2750
$ws.Range("K58").Value = 1638.0741
$ws.Range("L58").Value = 2750
$ws.Range("M58").Value = -1435.0741
$ws.Range("N58").Value = -3156
$ws.Range("H62").Value = 2322.2222
$ws.Range("I62").Value = 2266.6667
$ws.Range("K62").Value = 2266.6667
$ws.Range("M62").Value = -1642.6667
$ws.Range("H65").Value = 2322.2222
$ws.Range("I65").Value = 2266.6667
$ws.Range("K65").Value = 11333.3335
$ws.Range("M65").Value = -8213.333500000001
$ws.Range("H107").Value = 1463
$ws.Range("I107").Value = 1577.8235
$ws.Range("J107").Value = 975
$ws.Range("K107").Value = 1577.8235
$ws.Range("L107").Value = 975
$ws.Range("M107").Value = 342.1765
$ws.Range("N107").Value = -4815
$ws.Range("H132").Value = 1779.2069
$ws.Range("I132").Value = 1291.625
$ws.Range("J132").Value = 4119.6
$ws.Range("K132").Value = 3874.875
$ws.Range("L132").Value = 12358.8
$ws.Range("M132").Value = -1344.875
$ws.Range("N132").Value = -17418.8
$ws.Range("H136").Value = 1714.7587
$ws.Range("I136").Value = 1638.0741
$ws.Range("J136").Value = 2750
$ws.Range("K136").Value = 4914.2223
$ws.Range("L136").Value = 8250
$ws.Range("M136").Value = -2364.2223
$ws.Range("N136").Value = -13350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 40.576923
$ws.Range("I12").Value = 73.7
$ws.Range("J12").Value = 19.875
$ws.Range("K12").Value = 221.1
$ws.Range("L12").Value = 59.625
$ws.Range("M12").Value = -48.10000000000002
$ws.Range("N12").Value = -405.625
$ws.Range("H98").Value = 1429671.2
$ws.Range("I98").Value = 1100
$ws.Range("J98").Value = 2001099.8
$ws.Range("K98").Value = 3300
$ws.Range("L98").Value = 6003299.4
$ws.Range("M98").Value = -1802
$ws.Range("N98").Value = -6006295.4
$ws.Range("H122").Value = 1705.6428
$ws.Range("I122").Value = 2250
$ws.Range("J122").Value = 1557.1818
$ws.Range("K122").Value = 20250
$ws.Range("L122").Value = 14014.6362
$ws.Range("M122").Value = -17800
$ws.Range("N122").Value = -18914.6362
$ws.Range("H127").Value = 717.3
$ws.Range("J127").Value = 717.3
$ws.Range("L127").Value = 2151.9
$ws.Range("N127").Value = -12071.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 393408.1
$ws.Range("I46").Value = 5240
$ws.Range("J46").Value = 490450.12
$ws.Range("K46").Value = 5240
$ws.Range("L46").Value = 490450.12
$ws.Range("M46").Value = -5052
$ws.Range("N46").Value = -490826.12
$ws.Range("H48").Value = 19996
$ws.Range("I48").Value = 4990
$ws.Range("J48").Value = 30000
$ws.Range("K48").Value = 4990
$ws.Range("L48").Value = 30000
$ws.Range("M48").Value = -4329
$ws.Range("N48").Value = -31322
$ws.Range("H55").Value = 553.8333
$ws.Range("I55").Value = 608.7
$ws.Range("J55").Value = 485.25
$ws.Range("K55").Value = 608.7
$ws.Range("L55").Value = 485.25
$ws.Range("M55").Value = -435.7
$ws.Range("N55").Value = -831.25
$ws.Range("H61").Value = 1700.75
$ws.Range("I61").Value = 1831.75
$ws.Range("J61").Value = 1438.75
$ws.Range("K61").Value = 1831.75
$ws.Range("L61").Value = 1438.75
$ws.Range("M61").Value = -1629.75
$ws.Range("N61").Value = -1842.75
$ws.Range("H113").Value = 1700.75
$ws.Range("I113").Value = 1831.75
$ws.Range("J113").Value = 1438.75
$ws.Range("K113").Value = 1831.75
$ws.Range("L113").Value = 1438.75
$ws.Range("M113").Value = 338.25
$ws.Range("N113").Value = -5778.75
$ws.Range("H132").Value = 4417.9507
$ws.Range("I132").Value = 3086.1143
$ws.Range("J132").Value = 6210.8076
$ws.Range("K132").Value = 9258.3429
$ws.Range("L132").Value = 18632.4228
$ws.Range("M132").Value = -6728.3429
$ws.Range("N132").Value = -23692.4228
$ws.Range("H136").Value = 2627.9827
$ws.Range("I136").Value = 1932.5143
$ws.Range("J136").Value = 3686.3044
$ws.Range("K136").Value = 5797.5429
$ws.Range("L136").Value = 11058.9132
$ws.Range("M136").Value = -3247.5429
$ws.Range("N136").Value = -16158.9132

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 572
$ws.Range("I113").Value = 506.85715
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 1520.57145
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = 649.4285500000001
$ws.Range("N113").Value = -6740
$ws.Range("H132").Value = 1241.7715
$ws.Range("I132").Value = 934.25
$ws.Range("J132").Value = 1500.7368
$ws.Range("K132").Value = 2802.75
$ws.Range("L132").Value = 4502.2104
$ws.Range("M132").Value = -272.75
$ws.Range("N132").Value = -9562.2104
$ws.Range("H136").Value = 816.19446
$ws.Range("I136").Value = 608.26086
$ws.Range("J136").Value = 1184.0769
$ws.Range("K136").Value = 1824.78258
$ws.Range("L136").Value = 3552.2307
$ws.Range("M136").Value = 725.2174199999999
$ws.Range("N136").Value = -8652.2307

